$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking values
# like "1.0000" or "301.37" keep their exact text representation instead
# of being auto-converted to a number.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '23.042.96'
$ws.Range('E2').Value = '  -3.38%  '
$ws.Range('D3').Value = '1.600.91'
$ws.Range('E3').Value = '  -2.43%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '1.0000'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').Value = '301.37'
$ws.Range('D7').Value = '0.3782'
$ws.Range('E7').Value = '  -2.48%  '
$ws.Range('D8').Value = '0.3638'
$ws.Range('E8').Value = '  -4.96%  '
$ws.Range('D9').Value = '50.09'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').Value = '1.262'
$ws.Range('E10').Value = '  -4.78%  '
$ws.Range('D11').Value = '1.0000'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '0.08124'
$ws.Range('E12').Value = '  -3.15%  '
$ws.Range('D13').Value = '22.63'
$ws.Range('E13').Value = '  -4.42%  '
$ws.Range('D14').Value = '6.604'
$ws.Range('E14').Value = '  -5.27%  '
$ws.Range('E15').Value = '  -6.45%  '
$ws.Range('E16').Value = '  -5.00%  '
$ws.Range('D17').Value = '1.607.27'
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('D18').Value = '91.75'
$ws.Range('E18').Value = '  -2.01%  '
$ws.Range('D19').Value = '0.06821'
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('E20').Value = '  -6.04%  '
$ws.Range('D21').Value = '6.561'
$ws.Range('E21').Value = '  -4.51%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').Value = '13.09'
$ws.Range('E23').Value = '  -3.61%  '
$ws.Range('B24').Value = 'WrappedBTC'
$ws.Range('C24').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D24').Value = '23.036.70'
$ws.Range('E24').Value = '  -3.43%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.361'
$ws.Range('E25').Value = '  -3.34%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '2.846'
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '21.04'
$ws.Range('E27').Value = '  -3.53%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '150.40'
$ws.Range('E28').Value = '  -1.73%  '
$ws.Range('B29').Value = 'HuobiToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D29').Value = '5.254'
$ws.Range('E29').Value = '  -5.70%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '133.89'
$ws.Range('E30').Value = '  -2.01%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '6.847'
$ws.Range('E31').Value = '  -10.34%  '
$ws.Range('B32').Value = 'WEMIXTOKEN'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = '2.237'
$ws.Range('E32').Value = '  -10.34%  '
$ws.Range('B33').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C33').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D33').Value = '1.787.32'
$ws.Range('E33').Value = '  -1.81%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '0.9667'
$ws.Range('E34').Value = '  -1.18%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.07592'
$ws.Range('E35').Value = '  -5.00%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '10.39'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '6.283'
$ws.Range('E37').Value = '  -4.36%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.02706'
$ws.Range('E38').Value = '  -6.52%  '
$ws.Range('E39').Value = '  -4.36%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.08886'
$ws.Range('E40').Value = '  -2.09%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.373'
$ws.Range('E41').Value = '  -2.81%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.7029'
$ws.Range('E42').Value = '  -6.09%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '12.46'
$ws.Range('E43').Value = '  -6.12%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '15.26'
$ws.Range('E44').Value = '  -7.41%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.6632'
$ws.Range('E45').Value = '  -3.64%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '2.310'
$ws.Range('E46').Value = '  -4.39%  '
$ws.Range('D47').Value = '0.9986'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = '3.991'
$ws.Range('E48').Value = '  -2.29%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '132.47'
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('B50').Value = 'Flow'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D50').Value = '1.237'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('D51').Value = '0.07920'
$ws.Range('E51').Value = '  -3.66%  '

# Restore the default style on the price column so only the inline text
# values change (matches the original unstyled cells).
$ws.Range("D2:D51").Style = "Normal"

